# Add a new field "refrigerator_text" (type "text") to the "model" sheet so
# the form knows to show a popup/text message when the model record is
# invalid and the user attempts to create a new instance.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")
$ws.Activate()

$ws.Cells.Item(7, 1).Value = "refrigerator_text"
$ws.Cells.Item(7, 2).Value = "text"

$ws.Rows.Item(7).RowHeight = 12.75

$ws.Range("B7").Select()
